$wb = $excel.ActiveWorkbook

# ALC!row51
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 1334.909
$ws.Range("I51").Value = 900
$ws.Range("J51").Value = 1378.4
$ws.Range("K51").Value = 900
$ws.Range("L51").Value = 1378.4
$ws.Range("M51").Value = -416
$ws.Range("N51").Value = -2346.4

# ALC!row86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 3684.7222
$ws.Range("I86").Value = 3669.6365
$ws.Range("J86").Value = 3708.4285
$ws.Range("K86").Value = 3669.6365
$ws.Range("L86").Value = 3708.4285
$ws.Range("M86").Value = -2546.6365
$ws.Range("N86").Value = -5954.4285

# ALC!row89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 3684.7222
$ws.Range("I89").Value = 3669.6365
$ws.Range("J89").Value = 3708.4285
$ws.Range("K89").Value = 18348.1825
$ws.Range("L89").Value = 18542.1425
$ws.Range("M89").Value = -12732.1825
$ws.Range("N89").Value = -29774.1425

# ALC!row101
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 1423.1111
$ws.Range("I101").Value = 1490.375
$ws.Range("K101").Value = 4471.125
$ws.Range("M101").Value = -2849.125

# ALC!row125
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 2000.7142
$ws.Range("I125").Value = 1266
$ws.Range("J125").Value = 2123.1667
$ws.Range("K125").Value = 11394
$ws.Range("L125").Value = 19108.5003
$ws.Range("M125").Value = -8934
$ws.Range("N125").Value = -24028.5003

# ARM!row61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2398.5
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()

# ARM!row74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 822.1579
$ws.Range("I74").Value = 773.1875
$ws.Range("J74").Value = 1083.3334
$ws.Range("K74").Value = 773.1875
$ws.Range("L74").Value = 1083.3334
$ws.Range("M74").Value = 100.8125
$ws.Range("N74").Value = -2831.3334

# ARM!row77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 822.1579
$ws.Range("I77").Value = 773.1875
$ws.Range("J77").Value = 1083.3334
$ws.Range("K77").Value = 3865.9375
$ws.Range("L77").Value = 5416.666999999999
$ws.Range("M77").Value = 502.0625
$ws.Range("N77").Value = -14152.667

# ARM!row122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1618.1333
$ws.Range("I122").Value = 2128.6667
$ws.Range("J122").Value = 1277.7778
$ws.Range("K122").Value = 6386.000100000001
$ws.Range("L122").Value = 3833.3334
$ws.Range("M122").Value = -3936.000100000001
$ws.Range("N122").Value = -8733.3334

# ARM!row132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2582.8445
$ws.Range("I132").Value = 1689.9678
$ws.Range("J132").Value = 4559.9287
$ws.Range("K132").Value = 5069.903399999999
$ws.Range("L132").Value = 13679.7861
$ws.Range("M132").Value = -2539.903399999999
$ws.Range("N132").Value = -18739.7861

# ARM!row136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2398.5
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

# BSM!row134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1457.4
$ws.Range("I134").Value = 1419.3334
$ws.Range("J134").Value = 1800
$ws.Range("K134").Value = 4258.0002
$ws.Range("L134").Value = 5400
$ws.Range("M134").Value = -1723.0002
$ws.Range("N134").Value = -10470

# CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2125.164
$ws.Range("I31").Value = 1534.8684
$ws.Range("J31").Value = 2898.6553
$ws.Range("K31").Value = 1534.8684
$ws.Range("L31").Value = 2898.6553
$ws.Range("M31").Value = -1239.8684
$ws.Range("N31").Value = -3488.6553

# CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2125.164
$ws.Range("I34").Value = 1534.8684
$ws.Range("J34").Value = 2898.6553
$ws.Range("K34").Value = 1534.8684
$ws.Range("L34").Value = 2898.6553
$ws.Range("M34").Value = -1332.8684
$ws.Range("N34").Value = -3302.6553

# CRP!row58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1740.3226
$ws.Range("I58").Value = 1653.4482
$ws.Range("J58").Value = 3000
$ws.Range("K58").Value = 1653.4482
$ws.Range("L58").Value = 3000
$ws.Range("M58").Value = -1450.4482
$ws.Range("N58").Value = -3406

# CRP!row136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1740.3226
$ws.Range("I136").Value = 1653.4482
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 4960.3446
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -2410.3446
$ws.Range("N136").Value = -14100

# CUL!row51
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 5233.3335

# CUL!row110
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 8607.143
$ws.Range("I110").Value = 6000
$ws.Range("J110").Value = 9650
$ws.Range("K110").Value = 18000
$ws.Range("L110").Value = 28950
$ws.Range("M110").Value = -13910
$ws.Range("N110").Value = -37130

# CUL!row131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 922.4253
$ws.Range("I131").Value = 655.75
$ws.Range("J131").Value = 965.0933
$ws.Range("K131").Value = 1967.25
$ws.Range("L131").Value = 2895.2799
$ws.Range("M131").Value = 3072.75
$ws.Range("N131").Value = -12975.2799

# GSM!row80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3105.3333
$ws.Range("I80").Value = 2722.5
$ws.Range("J80").Value = 3542.8572
$ws.Range("K80").Value = 2722.5
$ws.Range("L80").Value = 3542.8572
$ws.Range("M80").Value = -1724.5
$ws.Range("N80").Value = -5538.8572

# GSM!row83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3105.3333
$ws.Range("I83").Value = 2722.5
$ws.Range("J83").Value = 3542.8572
$ws.Range("K83").Value = 13612.5
$ws.Range("L83").Value = 17714.286
$ws.Range("M83").Value = -8620.5
$ws.Range("N83").Value = -27698.286

# GSM!row105
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H105").Value = 27471.75
$ws.Range("J105").Value = 27471.75
$ws.Range("L105").Value = 27471.75
$ws.Range("N105").Value = -34459.75

# LTW!row7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5001559
$ws.Range("I7").Value = 6667900.5
$ws.Range("K7").Value = 6667900.5
$ws.Range("M7").Value = -6667788.5

# LTW!row122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 9808031
$ws.Range("I122").Value = 83334780
$ws.Range("J122").Value = 4464.067
$ws.Range("K122").Value = 250004340
$ws.Range("L122").Value = 13392.201
$ws.Range("M122").Value = -250001890
$ws.Range("N122").Value = -18292.201

# LTW!row126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 5001559
$ws.Range("I126").Value = 6667900.5
$ws.Range("K126").Value = 20003701.5
$ws.Range("M126").Value = -20001231.5

# LTW!row132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4160.154
$ws.Range("I132").Value = 4494.4
$ws.Range("K132").Value = 13483.2
$ws.Range("M132").Value = -10953.2

# WVR!row100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 477204.66
$ws.Range("I100").Value = 1085.8334
$ws.Range("K100").Value = 2171.6668
$ws.Range("M100").Value = -1630.6668

# WVR!row132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2735.1853
$ws.Range("I132").Value = 2794.6
$ws.Range("J132").Value = 2565.4285
$ws.Range("K132").Value = 8383.799999999999
$ws.Range("L132").Value = 7696.2855
$ws.Range("M132").Value = -5853.799999999999
$ws.Range("N132").Value = -12756.2855

Write-Output "edits applied"